# Rename the FV2310/"old" and FV2404/"new" header columns, (re)build the
# Excel Table (ListObject) over the data range and freeze the header row,
# matching the "adapt column header formatting to respective input file
# names" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- 1. Rename the header row -------------------------------------------------
# Left half of the diff table used the "_old" suffix, now "_FV2310".
$ws.Range("A1").Value = "Segmentname_FV2310"
$ws.Range("B1").Value = "Segmentgruppe_FV2310"
$ws.Range("C1").Value = "Segment_FV2310"
$ws.Range("D1").Value = "Datenelement_FV2310"
$ws.Range("E1").Value = "Segment ID_FV2310"
$ws.Range("F1").Value = "Code_FV2310"
$ws.Range("G1").Value = "Qualifier_FV2310"
$ws.Range("H1").Value = "Beschreibung_FV2310"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("J1").Value = "Bedingung_FV2310"

# Middle "diff" column is unchanged.
$ws.Range("K1").Value = "diff"

# Right half used the "_new" suffix, now "_FV2404".
$ws.Range("L1").Value = "Segmentname_FV2404"
$ws.Range("M1").Value = "Segmentgruppe_FV2404"
$ws.Range("N1").Value = "Segment_FV2404"
$ws.Range("O1").Value = "Datenelement_FV2404"
$ws.Range("P1").Value = "Segment ID_FV2404"
$ws.Range("Q1").Value = "Code_FV2404"
$ws.Range("R1").Value = "Qualifier_FV2404"
$ws.Range("S1").Value = "Beschreibung_FV2404"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("U1").Value = "Bedingung_FV2404"

# --- 2. Turn the used range into an Excel Table (ListObject) ------------------
$dataRange = $ws.Range("A1:U93")
$table = $ws.ListObjects.Add(1, $dataRange, [System.Reflection.Missing]::Value, 1)
$table.Name = "Table1"

# --- 3. Freeze the header row --------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
